$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = "scott"
$ws.Range("B2").Value = "son"
$ws.Range("C2").Value = "1,23,4,5"

# Delete rows 3 to 6 (which currently contain data)
$ws.Range("A3:C6").EntireRow.Delete()
